$wb = $excel.ActiveWorkbook

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 4491.154
$ws.Range("I39").Value = 1198.25
$ws.Range("J39").Value = 9759.799999999999
$ws.Range("K39").Value = 3594.75
$ws.Range("L39").Value = 29279.4
$ws.Range("M39").Value = -3298.75
$ws.Range("N39").Value = -29871.4

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 11006
$ws.Range("I129").Value = 1151
$ws.Range("K129").Value = 3453
$ws.Range("M129").Value = 1547

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2099.4167
$ws.Range("I137").Value = 1844
$ws.Range("J137").Value = 2384.8823
$ws.Range("K137").Value = 5532
$ws.Range("L137").Value = 7154.646900000001
$ws.Range("M137").Value = -2982
$ws.Range("N137").Value = -12254.6469

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 590114.2
$ws.Range("I2").Value = 866703.8
$ws.Range("K2").Value = 866703.8
$ws.Range("M2").Value = -866590.8

# ARM row 19
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1332.8334
$ws.Range("I45").Value = 1399.5
$ws.Range("J45").Value = 1199.5
$ws.Range("K45").Value = 1399.5
$ws.Range("L45").Value = 1199.5
$ws.Range("M45").Value = -1022.5
$ws.Range("N45").Value = -1953.5

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 333334660
$ws.Range("I61").Value = 333334660
$ws.Range("K61").Value = 333334660
$ws.Range("M61").Value = -333334448

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3979.3333
$ws.Range("I63").Value = 3979.3333
$ws.Range("K63").Value = 3979.3333
$ws.Range("M63").Value = -3293.3333

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3979.3333
$ws.Range("I66").Value = 3979.3333
$ws.Range("K66").Value = 19896.6665
$ws.Range("M66").Value = -16464.6665

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 143829.86
$ws.Range("J110").Value = 1250
$ws.Range("L110").Value = 1250
$ws.Range("N110").Value = -5340

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 590114.2
$ws.Range("I116").Value = 866703.8
$ws.Range("K116").Value = 866703.8
$ws.Range("M116").Value = -864409.8

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3812.7585
$ws.Range("I122").Value = 3039.6296
$ws.Range("K122").Value = 9118.888800000001
$ws.Range("M122").Value = -6668.888800000001

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5940116.5
$ws.Range("I132").Value = 3848908
$ws.Range("J132").Value = 15002018
$ws.Range("K132").Value = 11546724
$ws.Range("L132").Value = 45006054
$ws.Range("M132").Value = -11544194
$ws.Range("N132").Value = -45011114

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 333334660
$ws.Range("I136").Value = 333334660
$ws.Range("K136").Value = 1000003980
$ws.Range("M136").Value = -1000001430

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 590114.2
$ws.Range("I3").Value = 866703.8
$ws.Range("K3").Value = 866703.8
$ws.Range("M3").Value = -866589.8

# BSM row 30
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 39999
$ws.Range("J30").Value = 39999
$ws.Range("L30").Value = 39999
$ws.Range("N30").Value = -40249

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 113004.336
$ws.Range("I107").Value = 2130
$ws.Range("K107").Value = 2130
$ws.Range("M107").Value = -210

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25757200
$ws.Range("I134").Value = 28618722
$ws.Range("J134").Value = 3498.5
$ws.Range("K134").Value = 85856166
$ws.Range("L134").Value = 10495.5
$ws.Range("M134").Value = -85853631
$ws.Range("N134").Value = -15565.5

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1360285.5
$ws.Range("I16").Value = 1553897.8
$ws.Range("K16").Value = 1553897.8
$ws.Range("M16").Value = -1553610.8

# CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1053.5
$ws.Range("J19").Value = 1700
$ws.Range("L19").Value = 1700
$ws.Range("N19").Value = -2040

# CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 1053.5
$ws.Range("J24").Value = 1700
$ws.Range("L24").Value = 1700
$ws.Range("N24").Value = -2040

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2785.0925
$ws.Range("I31").Value = 2045.6945
$ws.Range("J31").Value = 4263.8887
$ws.Range("K31").Value = 2045.6945
$ws.Range("L31").Value = 4263.8887
$ws.Range("M31").Value = -1750.6945
$ws.Range("N31").Value = -4853.8887

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2785.0925
$ws.Range("I34").Value = 2045.6945
$ws.Range("J34").Value = 4263.8887
$ws.Range("K34").Value = 2045.6945
$ws.Range("L34").Value = 4263.8887
$ws.Range("M34").Value = -1843.6945
$ws.Range("N34").Value = -4667.8887

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3485.9375
$ws.Range("I99").Value = 3345.8333
$ws.Range("K99").Value = 3345.8333
$ws.Range("M99").Value = -1847.8333

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1360285.5
$ws.Range("I113").Value = 1553897.8
$ws.Range("K113").Value = 1553897.8
$ws.Range("M113").Value = -1551727.8

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3242.4614
$ws.Range("I122").Value = 3179.3333
$ws.Range("K122").Value = 9537.999899999999
$ws.Range("M122").Value = -7087.999899999999

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3485.9375
$ws.Range("I126").Value = 3345.8333
$ws.Range("K126").Value = 10037.4999
$ws.Range("M126").Value = -7567.499899999999

# CUL row 128
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 29277.6
$ws.Range("I18").Value = 19199.5
$ws.Range("J18").Value = 35996.332
$ws.Range("K18").Value = 19199.5
$ws.Range("L18").Value = 35996.332
$ws.Range("M18").Value = -18906.5
$ws.Range("N18").Value = -36582.332

# GSM row 21
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 32856.57
$ws.Range("I21").Value = 31922.46
$ws.Range("K21").Value = 31922.46
$ws.Range("M21").Value = -31749.46

# GSM row 30
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 32856.57
$ws.Range("I30").Value = 31922.46
$ws.Range("K30").Value = 31922.46
$ws.Range("M30").Value = -31817.46

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 15856.857
$ws.Range("I57").Value = 12666.5
$ws.Range("K57").Value = 12666.5
$ws.Range("M57").Value = -11846.5

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1860.6818
$ws.Range("I97").Value = 1863.8235
$ws.Range("K97").Value = 1863.8235
$ws.Range("M97").Value = -1367.8235

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3511.1667
$ws.Range("I126").Value = 3615.6
$ws.Range("K126").Value = 10846.8
$ws.Range("M126").Value = -8376.799999999999

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1766163.1
$ws.Range("I132").Value = 2238697.5
$ws.Range("J132").Value = 2035.4667
$ws.Range("K132").Value = 6716092.5
$ws.Range("L132").Value = 6106.4001
$ws.Range("M132").Value = -6713562.5
$ws.Range("N132").Value = -11166.4001

# LTW row 23
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 7745.75
$ws.Range("I23").Value = 4992
$ws.Range("K23").Value = 4992
$ws.Range("M23").Value = -4762

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6468
$ws.Range("I61").Value = 6853.25
$ws.Range("J61").Value = 1845
$ws.Range("K61").Value = 6853.25
$ws.Range("L61").Value = 1845
$ws.Range("M61").Value = -6651.25
$ws.Range("N61").Value = -2249

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2348.5
$ws.Range("I93").Value = 1417.6
$ws.Range("J93").Value = 3900
$ws.Range("K93").Value = 1417.6
$ws.Range("L93").Value = 3900
$ws.Range("M93").Value = -169.5999999999999
$ws.Range("N93").Value = -6396

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6468
$ws.Range("I113").Value = 6853.25
$ws.Range("J113").Value = 1845
$ws.Range("K113").Value = 6853.25
$ws.Range("L113").Value = 1845
$ws.Range("M113").Value = -4683.25
$ws.Range("N113").Value = -6185

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3231.6
$ws.Range("J122").Value = 3500
$ws.Range("L122").Value = 10500
$ws.Range("N122").Value = -15400

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16676188
$ws.Range("I132").Value = 17867238
$ws.Range("K132").Value = 53601714
$ws.Range("M132").Value = -53599184

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1986.9796
$ws.Range("I136").Value = 1944
$ws.Range("J136").Value = 2007.8182
$ws.Range("K136").Value = 5832
$ws.Range("L136").Value = 6023.4546
$ws.Range("M136").Value = -3282
$ws.Range("N136").Value = -11123.4546

# WVR row 87
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 59999
$ws.Range("J87").Value = 59999
$ws.Range("L87").Value = 59999
$ws.Range("N87").Value = -62495

# WVR row 90
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H90").Value = 59999
$ws.Range("J90").Value = 59999
$ws.Range("L90").Value = 179997
$ws.Range("N90").Value = -192477

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2365.7778
$ws.Range("I100").Value = 3132.1667
$ws.Range("K100").Value = 6264.3334
$ws.Range("M100").Value = -5723.3334

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21741580
$ws.Range("I132").Value = 33335422
$ws.Range("J132").Value = 3129.125
$ws.Range("K132").Value = 100006266
$ws.Range("L132").Value = 9387.375
$ws.Range("M132").Value = -100003736
$ws.Range("N132").Value = -14447.375

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 27780746
$ws.Range("I136").Value = 31253076
$ws.Range("K136").Value = 93759228
$ws.Range("M136").Value = -93756678
